$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$data = @(
    @(2, 'Última actualización: 11:17:08', $null, $null, $null, $null),
    @(3, 'Total filas: 194', $null, $null, $null, $null),
    @(89, '08:19:33', '09:02', '23_HERNANDEZ', 43, 'LP1912'),
    @(90, '08:48:09', '09:02', '215A_EL PATO', 14, 'LP1912'),
    @(156, '11:17:08', '11:17', '86_EST CHICA-ESC AGRARIA', 0, 'LP1912'),
    @(157, '11:17:08', '11:17', '16_SANTA ANA', 0, 'LP1912'),
    @(158, '11:17:08', '11:18', '15_ABASTO', 1, 'LP1912'),
    @(159, '09:25:56', '11:19', '86_EST CHICA-ESC AGRARIA', 114, 'LP1912'),
    @(160, '09:25:56', '11:20', '26_HERNANDEZ', 115, 'LP1912'),
    @(161, '10:12:35', '11:21', '26_HERNANDEZ', 69, 'LP1912'),
    @(162, '10:52:48', '11:26', '225_C ROCA-H SUR', 34, 'LP1912'),
    @(163, '10:12:35', '11:27', '225_C ROCA-H SUR', 75, 'LP1912'),
    @(164, '10:12:35', '11:32', '81_EL PELIGRO', 80, 'LP1912'),
    @(165, '10:52:48', '11:34', '23_HERNANDEZ', 42, 'LP1912'),
    @(166, '11:17:08', '11:35', '23_HERNANDEZ', 18, 'LP1912'),
    @(167, '10:52:48', '11:41', '17_ROMERO', 49, 'LP1912'),
    @(168, '10:12:35', '11:42', '17_ROMERO', 90, 'LP1912'),
    @(169, '10:52:48', '11:43', '10_OLMOS', 51, 'LP1912'),
    @(170, '10:12:35', '11:51', '215B_EL PATO', 99, 'LP1912'),
    @(171, '10:12:35', '11:56', '10_OLMOS', 104, 'LP1912'),
    @(172, '10:52:48', '11:58', '225_GOMEZ', 66, 'LP1912'),
    @(173, '10:12:35', '11:58', '16_P MOR-167 Y 521', 106, 'LP1912'),
    @(174, '10:12:35', '11:59', '225_GOMEZ', 107, 'LP1912'),
    @(175, '10:12:35', '12:04', '84_COLONIA URQUIZA-ESC 49', 112, 'LP1912'),
    @(176, '11:17:08', '12:05', '23_HERNANDEZ', 48, 'LP1912'),
    @(177, '10:12:35', '12:06', '16_P MOR-SANTA ANA', 114, 'LP1912'),
    @(178, '10:12:35', '12:06', '14_ABASTO', 114, 'LP1912'),
    @(179, '11:17:08', '12:07', '14_ABASTO', 50, 'LP1912'),
    @(180, '11:17:08', '12:07', '16_P MOR-SANTA ANA', 50, 'LP1912'),
    @(181, '10:52:48', '12:09', '84_COLONIA URQUIZA-ESC 49', 77, 'LP1912'),
    @(182, '11:17:08', '12:09', '10_OLMOS', 52, 'LP1912'),
    @(183, '11:17:08', '12:11', '84_COLONIA URQUIZA-ESC 49', 54, 'LP1912'),
    @(184, '10:52:48', '12:18', '14_ABASTO', 86, 'LP1912'),
    @(185, '10:52:48', '12:20', '215A_EL PATO', 88, 'LP1912'),
    @(186, '10:52:48', '12:21', '26_HERNANDEZ', 89, 'LP1912'),
    @(187, '11:17:08', '12:21', '14_ABASTO', 64, 'LP1912'),
    @(188, '11:17:08', '12:21', '215A_EL PATO', 64, 'LP1912'),
    @(189, '10:52:48', '12:24', '17_ROMERO', 92, 'LP1912'),
    @(190, '10:52:48', '12:36', '27_EL RETIRO', 104, 'LP1912'),
    @(191, '11:17:08', '12:37', '27_EL RETIRO', 80, 'LP1912'),
    @(192, '10:52:48', '12:38', '17_179 Y 38', 106, 'LP1912'),
    @(193, '11:17:08', '12:41', '10_OLMOS', 84, 'LP1912'),
    @(194, '11:17:08', '12:48', '17_ROMERO', 91, 'LP1912'),
    @(195, '10:52:48', '12:48', '11_ETCHEVERRY', 116, 'LP1912'),
    @(196, '11:17:08', '12:49', '11_ETCHEVERRY', 92, 'LP1912'),
    @(197, '11:17:08', '13:07', '16_P MOR-SANTA ANA', 110, 'LP1912'),
    @(198, '11:17:08', '13:14', '215D_EL PATO', 117, 'LP1912'),
    @(199, '11:17:08', '13:16', '11_ETCHEVERRY', 119, 'LP1912')
)
foreach ($row in $data) {
    $r = $row[0]
    for ($c = 1; $c -le 5; $c++) {
        $val = $row[$c]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$data = @(
    @(2, 'Última actualización: 11:17:08', $null, $null, $null, $null),
    @(3, 'Total filas: 23', $null, $null, $null, $null),
    @(27, '11:17:08', '12:21', '215A_EL PATO', 64, 'LP1912'),
    @(28, '11:17:08', '13:14', '215D_EL PATO', 117, 'LP1912')
)
foreach ($row in $data) {
    $r = $row[0]
    for ($c = 1; $c -le 5; $c++) {
        $val = $row[$c]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$data = @(
    @(2, 'Última actualización: 11:17:08', $null, $null, $null, $null),
    @(3, 'Total filas: 31', $null, $null, $null, $null),
    @(36, '11:17:08', '12:54', '215C_LA PLATA', 97, 'L6203')
)
foreach ($row in $data) {
    $r = $row[0]
    for ($c = 1; $c -le 5; $c++) {
        $val = $row[$c]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}
